$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 17.96285899517177
$ws.Range("C2").Value = 8.714392610023912
$ws.Range("D2").Value = 8.018325215040361
$ws.Range("E2").Value = 12.81274015110884
$ws.Range("F2").Value = 36.38381813520184
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("J2").Value = 10.0589417986614
$ws.Range("L2").Value = 10.95819561798788
$ws.Range("M2").Value = 16.91670573382476
$ws.Range("N2").Value = 19.49425914206329
$ws.Range("O2").Value = 28.01996576738952

$ws.Range("B3").Value = 17.62360461244322
$ws.Range("C3").Value = 8.516187486938055
$ws.Range("D3").Value = 8.023061227770759
$ws.Range("E3").Value = 12.84045685644295
$ws.Range("F3").Value = 36.45166563294431
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("J3").Value = 10.07565284077258
$ws.Range("L3").Value = 10.96114853881366
$ws.Range("M3").Value = 16.84883679745128
$ws.Range("N3").Value = 19.55327155283014
$ws.Range("O3").Value = 28.0866397459683

$ws.Range("B4").Value = 17.4149198962931
$ws.Range("C4").Value = 8.391035433927787
$ws.Range("D4").Value = 8.026676513517348
$ws.Range("E4").Value = 12.85848449738514
$ws.Range("F4").Value = 36.5014064122886
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("J4").Value = 10.08643842088905
$ws.Range("L4").Value = 10.96409290700709
$ws.Range("M4").Value = 16.80930443695769
$ws.Range("N4").Value = 19.59135150766554
$ws.Range("O4").Value = 28.13351603022223

$ws.Range("B5").Value = 17.32989696240417
$ws.Range("C5").Value = 8.339211677544647
$ws.Range("D5").Value = 8.028328138739599
$ws.Range("E5").Value = 12.8660853536798
$ws.Range("F5").Value = 36.52370500449146
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("J5").Value = 10.0909660190077
$ws.Range("L5").Value = 10.96557811437723
$ws.Range("M5").Value = 16.79374380987201
$ws.Range("N5").Value = 19.607334841653
$ws.Range("O5").Value = 28.15410887254716

$ws.Range("B6").Value = 17.31578340534568
$ws.Range("C6").Value = 8.33055802827554
$ws.Range("D6").Value = 8.028613176962969
$ws.Range("E6").Value = 12.86736285587677
$ws.Range("F6").Value = 36.52753009212138
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("J6").Value = 10.09172583047278
$ws.Range("L6").Value = 10.96584199153632
$ws.Range("M6").Value = 16.79119346037598
$ws.Range("N6").Value = 19.61001700996806
$ws.Range("O6").Value = 28.15761822581533

$ws.Range("B7").Value = 17.4137730235746
$ws.Range("C7").Value = 8.390339792512924
$ws.Range("D7").Value = 8.026698065060598
$ws.Range("E7").Value = 12.85858597413648
$ws.Range("F7").Value = 36.5016989293154
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("J7").Value = 10.08649894511634
$ws.Range("L7").Value = 10.96411178051006
$ws.Range("M7").Value = 16.80909234319093
$ws.Range("N7").Value = 19.5915651783074
$ws.Range("O7").Value = 28.13378772196575

$ws.Range("B8").Value = 17.84603063436433
$ws.Range("C8").Value = 8.646794423067995
$ws.Range("D8").Value = 8.019811683058027
$ws.Range("E8").Value = 12.82208772043336
$ws.Range("F8").Value = 36.40553254851267
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("J8").Value = 10.06459505047866
$ws.Range("L8").Value = 10.95897953747856
$ws.Range("M8").Value = 16.8928673387155
$ws.Range("N8").Value = 19.51422420345066
$ws.Range("O8").Value = 28.04172085213533

$ws.Range("B9").Value = 18.68580417460878
$ws.Range("C9").Value = 9.120460861147684
$ws.Range("D9").Value = 8.011898067653485
$ws.Range("E9").Value = 12.75849692486798
$ws.Range("F9").Value = 36.28121590091523
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("J9").Value = 10.02578884954097
$ws.Range("L9").Value = 10.95785121921187
$ws.Range("M9").Value = 17.07361231522706
$ws.Range("N9").Value = 19.37715044432065
$ws.Range("O9").Value = 27.9084127770012

$ws.Range("B10").Value = 19.29171971125373
$ws.Range("C10").Value = 9.448374094535202
$ws.Range("D10").Value = 8.009462304038788
$ws.Range("E10").Value = 12.71660455820034
$ws.Range("F10").Value = 36.2292168087523
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("J10").Value = 9.999781073023453
$ws.Range("L10").Value = 10.96241448552366
$ws.Range("M10").Value = 17.21572265118276
$ws.Range("N10").Value = 19.28525864185478
$ws.Range("O10").Value = 27.83941493266191

$ws.Range("B11").Value = 19.56372346085536
$ws.Range("C11").Value = 9.592765311697386
$ws.Range("D11").Value = 8.009081096547705
$ws.Range("E11").Value = 12.69858683440568
$ws.Range("F11").Value = 36.21412465905367
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("J11").Value = 9.988487638915771
$ws.Range("L11").Value = 10.96564845905674
$ws.Range("M11").Value = 17.2822325246727
$ws.Range("N11").Value = 19.2453523562526
$ws.Range("O11").Value = 27.8143368869849

$ws.Range("B12").Value = 19.66610556126149
$ws.Range("C12").Value = 9.646723727496893
$ws.Range("D12").Value = 8.009040651571013
$ws.Range("E12").Value = 12.69191283547247
$ws.Range("F12").Value = 36.20964161622678
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("J12").Value = 9.984288019403836
$ws.Range("L12").Value = 10.96703841236307
$ws.Range("M12").Value = 17.30767156349978
$ws.Range("N12").Value = 19.23051226962881
$ws.Range("O12").Value = 27.80574936777032

$ws.Range("B13").Value = 19.64408480880274
$ws.Range("C13").Value = 9.635135298810173
$ws.Range("D13").Value = 8.009044750607297
$ws.Range("E13").Value = 12.69334358560973
$ws.Range("F13").Value = 36.21055232008109
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("J13").Value = 9.985189065416224
$ws.Range("L13").Value = 10.96673172853469
$ws.Range("M13").Value = 17.30218178876153
$ws.Range("N13").Value = 19.23369628672435
$ws.Range("O13").Value = 27.80755839236495

$ws.Range("B14").Value = 19.57215946040286
$ws.Range("C14").Value = 9.597219096484078
$ws.Range("D14").Value = 8.00907569016184
$ws.Range("E14").Value = 12.69803477890013
$ws.Range("F14").Value = 36.21373114282498
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("J14").Value = 9.988140593327181
$ws.Range("L14").Value = 10.96575950532925
$ws.Range("M14").Value = 17.28432042096262
$ws.Range("N14").Value = 19.24412601692244
$ws.Range("O14").Value = 27.8136121604568

$ws.Range("B15").Value = 19.52801958121812
$ws.Range("C15").Value = 9.573899727104269
$ws.Range("D15").Value = 8.009108155274767
$ws.Range("E15").Value = 12.70092764861965
$ws.Range("F15").Value = 36.21583872077904
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("J15").Value = 9.989958501312991
$ws.Range("L15").Value = 10.96518548261167
$ws.Range("M15").Value = 17.27341234746301
$ws.Range("N15").Value = 19.25054986180282
$ws.Range("O15").Value = 27.81743869083292

$ws.Range("B16").Value = 19.27386209191121
$ws.Range("C16").Value = 9.438838816352527
$ws.Range("D16").Value = 8.009501791143467
$ws.Range("E16").Value = 12.71780292710943
$ws.Range("F16").Value = 36.23037549471332
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("J16").Value = 10.00052991389302
$ws.Range("L16").Value = 10.96222633928005
$ws.Range("M16").Value = 17.21141231557403
$ws.Range("N16").Value = 19.28790466876906
$ws.Range("O16").Value = 27.8411809568868

$ws.Range("B17").Value = 19.11694359877257
$ws.Range("C17").Value = 9.354735670428187
$ws.Range("D17").Value = 8.009929008040805
$ws.Range("E17").Value = 12.72842117959557
$ws.Range("F17").Value = 36.24148713787632
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("J17").Value = 10.00715258247365
$ws.Range("L17").Value = 10.96070677496988
$ws.Range("M17").Value = 17.1738441449857
$ws.Range("N17").Value = 19.31130542876007
$ws.Range("O17").Value = 27.85736337945688

$ws.Range("B18").Value = 19.02635204333016
$ws.Range("C18").Value = 9.305913779604406
$ws.Range("D18").Value = 8.010243182677138
$ws.Range("E18").Value = 12.73462637507747
$ws.Range("F18").Value = 36.24868420036381
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("J18").Value = 10.01101239178597
$ws.Range("L18").Value = 10.95994184624021
$ws.Range("M18").Value = 17.15241220314703
$ws.Range("N18").Value = 19.32494342608207
$ws.Range("O18").Value = 27.86726487181523

$ws.Range("B19").Value = 18.99562462771748
$ws.Range("C19").Value = 9.289307613175678
$ws.Range("D19").Value = 8.010361331449838
$ws.Range("E19").Value = 12.73674417105656
$ws.Range("F19").Value = 36.25125938932198
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("J19").Value = 10.01232796266098
$ws.Range("L19").Value = 10.95970162603826
$ws.Range("M19").Value = 17.14518643425181
$ws.Range("N19").Value = 19.32959170733915
$ws.Range("O19").Value = 27.87071927384534

$ws.Range("B20").Value = 19.13368336026623
$ws.Range("C20").Value = 9.363735201353499
$ws.Range("D20").Value = 8.009876450013937
$ws.Range("E20").Value = 12.72728072448044
$ws.Range("F20").Value = 36.24022086526027
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("J20").Value = 10.00644235110757
$ws.Range("L20").Value = 10.96085725399258
$ws.Range("M20").Value = 17.17782520247752
$ws.Range("N20").Value = 19.3087959115764
$ws.Range("O20").Value = 27.85557926439636

$ws.Range("B21").Value = 19.59330325229788
$ws.Range("C21").Value = 9.608375772799448
$ws.Range("D21").Value = 8.009063787329286
$ws.Range("E21").Value = 12.69665282421984
$ws.Range("F21").Value = 36.21276400597509
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("J21").Value = 9.987271572183634
$ws.Range("L21").Value = 10.9660405941819
$ws.Range("M21").Value = 17.28955998519712
$ws.Range("N21").Value = 19.24105519062005
$ws.Range("O21").Value = 27.81180934027643

$ws.Range("B22").Value = 19.89003723522547
$ws.Range("C22").Value = 9.764056180043747
$ws.Range("D22").Value = 8.009138034976843
$ws.Range("E22").Value = 12.67750351004437
$ws.Range("F22").Value = 36.20200043460336
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("J22").Value = 9.97519077390297
$ws.Range("L22").Value = 10.97039119768085
$ws.Range("M22").Value = 17.36405385631076
$ws.Range("N22").Value = 19.19836508535598
$ws.Range("O22").Value = 27.78850202362442

$ws.Range("B23").Value = 19.73202995252277
$ws.Range("C23").Value = 9.681361481444426
$ws.Range("D23").Value = 8.009043222692158
$ws.Range("E23").Value = 12.68764463068717
$ws.Range("F23").Value = 36.20708801009841
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("J23").Value = 9.981597608989613
$ws.Range("L23").Value = 10.9679815028669
$ws.Range("M23").Value = 17.32416564143827
$ws.Range("N23").Value = 19.22100515289262
$ws.Range("O23").Value = 27.80045625373973

$ws.Range("B24").Value = 19.12611648091806
$ws.Range("C24").Value = 9.359667974525857
$ws.Range("D24").Value = 8.009899997901774
$ws.Range("E24").Value = 12.72779601077949
$ws.Range("F24").Value = 36.24079082769713
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("J24").Value = 10.00676328361767
$ws.Range("L24").Value = 10.96078888378213
$ws.Range("M24").Value = 17.1760248466285
$ws.Range("N24").Value = 19.30992988921737
$ws.Range("O24").Value = 27.85638400018369

$ws.Range("B25").Value = 18.46012940756967
$ws.Range("C25").Value = 8.99569659200205
$ws.Range("D25").Value = 8.01344339678778
$ws.Range("E25").Value = 12.77484936061645
$ws.Range("F25").Value = 36.3079481096637
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("J25").Value = 10.0358455719295
$ws.Range("L25").Value = 10.95720529110227
$ws.Range("M25").Value = 17.02302582784204
$ws.Range("N25").Value = 19.41267864351315
$ws.Range("O25").Value = 27.93940292395234
